$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "At Work" (C) 1.0 entries need to move to other leave columns
# because they were actually Public Holidays / weekends leave entries.

# Move C -> E (Sick Leave)
$eRows = @(12, 13)
foreach ($r in $eRows) {
    $ws.Range("C$r").Value = ""
    $ws.Range("E$r").Value = "1.0"
}

# Move C -> F (Childcare Leave)
$fRows = @(20, 23, 24, 25)
foreach ($r in $fRows) {
    $ws.Range("C$r").Value = ""
    $ws.Range("F$r").Value = "1.0"
}

# Move C -> F and also set G (Annual Leave) for rows 26, 27
$fgRows = @(26, 27)
foreach ($r in $fgRows) {
    $ws.Range("C$r").Value = ""
    $ws.Range("F$r").Value = "1.0"
    $ws.Range("G$r").Value = "1.0"
}

# Move C -> G (Annual Leave)
$gRows = @(30, 31, 32, 33, 34, 37, 38, 41)
foreach ($r in $gRows) {
    $ws.Range("C$r").Value = ""
    $ws.Range("G$r").Value = "1.0"
}

# Update totals row 44
$ws.Range("C44").Value = "4.0"
$ws.Range("E44").Value = "2.0"
$ws.Range("F44").Value = "6.0"
$ws.Range("G44").Value = "10.0"

# Update signature date
$ws.Range("B50").Value = "07 - February - 2025"
